$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.857.13"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.741.75"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5126"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  +4.45%  "
$ws.Range("E9").Value = "  -5.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06085"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").Value = "1.738.35"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06983"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6329"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.498"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "25.876.24"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").Value = "1.958.26"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.080"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.488"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.096"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.506"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.813"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08262"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.613"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.395"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04399"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9681"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5974"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.672"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01545"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.904"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3813"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7278"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.877"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05470"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.238"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1099"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.478"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "